{"js": "// Remove the \"Ver no Jupiter ...\" line, the site footer/copyright line\n// (\"\u00a9 2020 . Contact: ...\") and the blank paragraph that separated them\n// from the preceding \"LOM3081: ...\" requirement line. This mirrors the\n// Jekyll site rebuild that dropped the page-chrome paragraphs from the\n// generated docx while leaving the rest of the \"Requisitos\" section and\n// the trailing blank / page-break paragraphs untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Ver no Jupiter ...\" and the \"\u00a9 ... Contact:\" paragraphs by\n// their text so the script is resilient to any surrounding content.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIdx === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIdx = i;\n  }\n  if (copyrightIdx === -1 && t.indexOf(\"Contact:\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx === -1 || copyrightIdx === -1) {\n  throw new Error(\"Could not locate the Jupiter/copyright paragraphs to remove.\");\n}\n\n// The blank paragraph right before the \"Ver no Jupiter ...\" line (only\n// when it is in fact empty) is removed together with it.\nlet blankIdx = -1;\nif (jupiterIdx - 1 >= 0 && items[jupiterIdx - 1].text === \"\") {\n  blankIdx = jupiterIdx - 1;\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nconst toDelete = [copyrightIdx, jupiterIdx];\nif (blankIdx !== -1) {\n  toDelete.push(blankIdx);\n}\ntoDelete.sort((a, b) => b - a);\n\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Ver no Jupiter ...\" line, the site footer/copyright line\n# (\"\u00a9 2020 . Contact: ...\") and the blank paragraph that separated them\n# from the preceding \"LOM3081: ...\" requirement line. This mirrors the\n# Jekyll site rebuild that dropped the page-chrome paragraphs from the\n# generated docx while leaving the rest of the \"Requisitos\" section and\n# the trailing blank / page-break paragraphs untouched.\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# Locate the \"Ver no Jupiter ...\" and the \"... Contact: ...\" (copyright)\n# paragraphs by their text so the script is resilient to any surrounding\n# content / exact paragraph numbering.\n$jupiterIndex = -1\n$copyrightIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIndex -eq -1 -and $text -like \"*Ver no Jupiter*\") {\n        $jupiterIndex = $i\n    }\n    if ($copyrightIndex -eq -1 -and $text -like \"*Contact:*\") {\n        $copyrightIndex = $i\n    }\n}\n\nif ($jupiterIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"Could not locate the Jupiter/copyright paragraphs to remove.\"\n}\n\n# The blank paragraph right before the \"Ver no Jupiter ...\" line (only\n# when it is in fact empty) is removed together with it.\n$blankIndex = -1\nif ($jupiterIndex -gt 1) {\n    $prevText = $d.Paragraphs.Item($jupiterIndex - 1).Range.Text\n    if ($prevText.Trim() -eq \"\") {\n        $blankIndex = $jupiterIndex - 1\n    }\n}\n\n# Delete from the bottom up so earlier indices stay valid.\n$d.Paragraphs.Item($copyrightIndex).Range.Delete()\n$d.Paragraphs.Item($jupiterIndex).Range.Delete()\nif ($blankIndex -ne -1) {\n    $d.Paragraphs.Item($blankIndex).Range.Delete()\n}\n"}
